# Apply the cryptocurrency price/volume update described by the commit.
# Values in column D (Price) must remain text (matching the source data,
# which stores numbers using "." as both thousands and decimal separators,
# e.g. "91.852.95"), so we force text storage via NumberFormat "@" and then
# restore each cell's original style to avoid leaving a residual number
# format on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D="91.852.95"; E="  -2.05%  " },
    @{ Row=3; D="3.326.25"; E="  -3.53%  " },
    @{ Row=4; E="  -0.02%  " },
    @{ Row=5; D="230.69"; E="  -2.53%  " },
    @{ Row=6; D="613.27"; E="  -4.00%  " },
    @{ Row=7; E="  -1.85%  " },
    @{ Row=8; E="  -2.93%  " },
    @{ Row=9; E="  +0.03%  " },
    @{ Row=10; D="0.956"; E="  +0.16%  " },
    @{ Row=11; D="3.325.85"; E="  -3.48%  " },
    @{ Row=12; D="42.71"; E="  +1.91%  " },
    @{ Row=13; E="  -1.56%  " },
    @{ Row=14; E="  +0.01%  " },
    @{ Row=15; D="91.721.29"; E="  -2.02%  " },
    @{ Row=16; D="3.945.91"; E="  -3.63%  " },
    @{ Row=17; E="  -2.85%  " },
    @{ Row=18; D="8.05"; E="  -3.01%  " },
    @{ Row=19; D="3.328.10"; E="  -3.57%  " },
    @{ Row=20; D="17.31"; E="  -1.65%  " },
    @{ Row=21; D="10.83"; E="  -3.80%  " },
    @{ Row=22; E="  +5.72%  " },
    @{ Row=23; D="491.51"; E="  -0.83%  " },
    @{ Row=24; D="0.440"; E="  -10.15%  " },
    @{ Row=25; D="6.49"; E="  -0.47%  " },
    @{ Row=26; E="  -4.87%  " },
    @{ Row=27; D="92.21"; E="  +1.77%  " },
    @{ Row=28; D="11.91"; E="  -0.32%  " },
    @{ Row=29; D="3.504.43"; E="  -3.51%  " },
    @{ Row=30; E="  -0.07%  " },
    @{ Row=31; E="  -4.70%  " },
    @{ Row=32; D="0.138"; E="  +2.05%  " },
    @{ Row=33; D="2.60"; E="  -4.69%  " },
    @{ Row=34; D="0.999"; E="  -0.12%  " },
    @{ Row=35; E="  -4.04%  " },
    @{ Row=36; D="28.28"; E="  -5.68%  " },
    @{ Row=37; D="0.527"; E="  -4.83%  " },
    @{ Row=38; D="566.21"; E="  +3.22%  " },
    @{ Row=39; B="RenderToken"; C="https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"; D="7.39"; E="  -2.49%  " },
    @{ Row=40; B="USDe"; C="https://coinranking.com/coin/exbfr2U-0+usde-usde"; D="1.00"; E="  -0.05%  " },
    @{ Row=41; D="0.149"; E="  -1.26%  " },
    @{ Row=42; E="  -4.89%  " },
    @{ Row=43; E="  -5.88%  " },
    @{ Row=44; D="23.69"; E="  -1.33%  " },
    @{ Row=45; D="1.68"; E="  -1.88%  " },
    @{ Row=46; E="  +0.62%  " },
    @{ Row=47; D="3.59"; E="  +1.13%  " },
    @{ Row=48; E="  -1.91%  " },
    @{ Row=49; D="2.10"; E="  -1.63%  " },
    @{ Row=50; B="Cosmos"; C="https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D="7.97"; E="  -0.18%  " },
    @{ Row=51; B="OKB"; C="https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; D="51.85"; E="  -2.19%  " }
)

foreach ($u in $updates) {
    $row = $u.Row
    foreach ($col in "B", "C", "D", "E") {
        if ($u.ContainsKey($col)) {
            $cell = $ws.Range("$col$row")
            if ($col -eq "D") {
                # Force text so numeric-looking strings (e.g. "230.69")
                # are not silently converted to a Double by Excel.
                $origStyle = $cell.Style
                $cell.NumberFormat = "@"
                $cell.Value = $u[$col]
                $cell.Style = $origStyle
            } else {
                $cell.Value = $u[$col]
            }
        }
    }
}

